# Auto-generated Excel COM-interop script applying the scheduled data refresh
# to the Zeromus_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates currentAveragePrice/LevePrice/LeveProfit columns (H:N) for the rows
# touched by this run; a few cells are added or removed where the source feed
# no longer reports (or newly reports) an HQ/NQ price for that leve.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2948.75
$ws.Range("I62").Value = 2511.4285
$ws.Range("J62").Value = 3561
$ws.Range("K62").Value = 2511.4285
$ws.Range("L62").Value = 3561
$ws.Range("M62").Value = -1887.4285
$ws.Range("N62").Value = -4809
# Row 65
$ws.Range("H65").Value = 2948.75
$ws.Range("I65").Value = 2511.4285
$ws.Range("J65").Value = 3561
$ws.Range("K65").Value = 12557.1425
$ws.Range("L65").Value = 17805
$ws.Range("M65").Value = -9437.1425
$ws.Range("N65").Value = -24045
# Row 98
$ws.Range("H98").Value = 1436.25
$ws.Range("I98").Value = 1127.5
$ws.Range("K98").Value = 1127.5
$ws.Range("M98").Value = 370.5
# Row 113
$ws.Range("H113").Value = 2754.75
$ws.Range("I113").Value = 1782
$ws.Range("J113").Value = 3196.9092
$ws.Range("K113").Value = 1782
$ws.Range("L113").Value = 3196.9092
$ws.Range("M113").Value = 1472
$ws.Range("N113").Value = -9704.9092
# Row 122
$ws.Range("H122").Value = 1436.25
$ws.Range("I122").Value = 1127.5
$ws.Range("K122").Value = 3382.5
$ws.Range("M122").Value = -932.5

$ws = $wb.Worksheets.Item("ARM")
# Row 52
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 30000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30636
# Row 61
$ws.Range("H61").Value = 2148.9285
$ws.Range("I61").Value = 2094.111
$ws.Range("J61").Value = 2247.6
$ws.Range("K61").Value = 2094.111
$ws.Range("L61").Value = 2247.6
$ws.Range("M61").Value = -1882.111
$ws.Range("N61").Value = -2671.6
# Row 102
$ws.Range("H102").Value = 1401.25
$ws.Range("I102").Value = 1401.25
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1401.25
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 220.75
$ws.Range("N102").ClearContents()
# Row 136
$ws.Range("H136").Value = 2148.9285
$ws.Range("I136").Value = 2094.111
$ws.Range("J136").Value = 2247.6
$ws.Range("K136").Value = 6282.333
$ws.Range("L136").Value = 6742.799999999999
$ws.Range("M136").Value = -3732.333
$ws.Range("N136").Value = -11842.8

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 8697774
$ws.Range("I86").Value = 10528502
$ws.Range("J86").Value = 1814.25
$ws.Range("K86").Value = 10528502
$ws.Range("L86").Value = 1814.25
$ws.Range("M86").Value = -10527379
$ws.Range("N86").Value = -4060.25
# Row 89
$ws.Range("H89").Value = 8697774
$ws.Range("I89").Value = 10528502
$ws.Range("J89").Value = 1814.25
$ws.Range("K89").Value = 52642510
$ws.Range("L89").Value = 9071.25
$ws.Range("M89").Value = -52636894
$ws.Range("N89").Value = -20303.25
# Row 94
$ws.Range("H94").Value = 4351.846
$ws.Range("I94").Value = 486.64706
$ws.Range("J94").Value = 11652.777
$ws.Range("K94").Value = 486.64706
$ws.Range("L94").Value = 11652.777
$ws.Range("M94").Value = -35.64706000000001
$ws.Range("N94").Value = -12554.777
# Row 134
$ws.Range("H134").Value = 1629.8096
$ws.Range("I134").Value = 1108
$ws.Range("J134").Value = 3299.6
$ws.Range("K134").Value = 3324
$ws.Range("L134").Value = 9898.799999999999
$ws.Range("M134").Value = -789
$ws.Range("N134").Value = -14968.8

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1760.3043
$ws.Range("I58").Value = 1039.909
$ws.Range("J58").Value = 2420.6667
$ws.Range("K58").Value = 1039.909
$ws.Range("L58").Value = 2420.6667
$ws.Range("M58").Value = -836.9090000000001
$ws.Range("N58").Value = -2826.6667
# Row 106
$ws.Range("H106").Value = 35333.332
$ws.Range("J106").Value = 35333.332
$ws.Range("L106").Value = 35333.332
$ws.Range("N106").Value = -37857.332
# Row 122
$ws.Range("H122").Value = 850
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2100
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 350
$ws.Range("N122").Value = -7900
# Row 136
$ws.Range("H136").Value = 1760.3043
$ws.Range("I136").Value = 1039.909
$ws.Range("J136").Value = 2420.6667
$ws.Range("K136").Value = 3119.727
$ws.Range("L136").Value = 7262.000100000001
$ws.Range("M136").Value = -569.7270000000003
$ws.Range("N136").Value = -12362.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 5183.737
$ws.Range("I3").Value = 3064
$ws.Range("J3").Value = 7539
$ws.Range("K3").Value = 9192
$ws.Range("L3").Value = 22617
$ws.Range("M3").Value = -9080
$ws.Range("N3").Value = -22841
# Row 120
$ws.Range("H120").Value = 5032.5
$ws.Range("I120").Value = 3376.6667
$ws.Range("J120").Value = 10000
$ws.Range("K120").Value = 10130.0001
$ws.Range("L120").Value = 30000
$ws.Range("M120").Value = -5292.000100000001
$ws.Range("N120").Value = -39676

$ws = $wb.Worksheets.Item("GSM")
# Row 23
$ws.Range("H23").Value = 3076.5454
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# Row 70
$ws.Range("H70").Value = 1733927.1
$ws.Range("I70").Value = 4904960.5
$ws.Range("J70").Value = 4272.727
$ws.Range("K70").Value = 4904960.5
$ws.Range("L70").Value = 4272.727
$ws.Range("M70").Value = -4904690.5
$ws.Range("N70").Value = -4812.727
# Row 73
$ws.Range("H73").Value = 1733927.1
$ws.Range("I73").Value = 4904960.5
$ws.Range("J73").Value = 4272.727
$ws.Range("K73").Value = 4904960.5
$ws.Range("L73").Value = 4272.727
$ws.Range("M73").Value = -4904024.5
$ws.Range("N73").Value = -6144.727
# Row 122
$ws.Range("H122").Value = 2558.8333
$ws.Range("I122").Value = 2558.8333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7676.499899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5226.499899999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1670.6342
$ws.Range("J7").Value = 1899.8
$ws.Range("L7").Value = 1899.8
$ws.Range("N7").Value = -2123.8
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
# Row 40
$ws.Range("H40").Value = 1989.909
$ws.Range("I40").Value = 1986.125
$ws.Range("K40").Value = 1986.125
$ws.Range("M40").Value = -1850.125
# Row 126
$ws.Range("H126").Value = 1670.6342
$ws.Range("J126").Value = 1899.8
$ws.Range("L126").Value = 5699.4
$ws.Range("N126").Value = -10639.4
# Row 136
$ws.Range("H136").Value = 7426.55
$ws.Range("I136").Value = 18103.715
$ws.Range("J136").Value = 1677.3077
$ws.Range("K136").Value = 54311.145
$ws.Range("L136").Value = 5031.9231
$ws.Range("M136").Value = -51761.145
$ws.Range("N136").Value = -10131.9231

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1653.871
$ws.Range("I136").Value = 1052.375
$ws.Range("J136").Value = 2295.4666
$ws.Range("K136").Value = 3157.125
$ws.Range("L136").Value = 6886.399800000001
$ws.Range("M136").Value = -607.125
$ws.Range("N136").Value = -11986.3998
